# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-06 07:13:35
#
# For every row in column G ("Recorded By") that contains a comma-separated
# list of recorders, rotate the list by moving the last entry to the front
# (e.g. "a, b, System" -> "System, a, b").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Count -gt 1) {
        $last = $parts[$parts.Count - 1]
        $rest = $parts[0..($parts.Count - 2)]
        $newParts = @($last) + $rest
        $newText = $newParts -join ", "
        $cell.Value = $newText
    }
}
